$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68. This shifts the existing rows 68-114
# down to 69-115, carrying their values/styles with them, and leaves a
# blank (but style-inheriting) row 68 to be populated below.
$ws.Rows.Item(68).Insert()

# Populate the new row 68 with the new data record.
$ws.Cells.Item(68, 1).Value = 10
$ws.Cells.Item(68, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(68, 3).Value = "La Araucanía"
$ws.Cells.Item(68, 4).Value = 44942
$ws.Cells.Item(68, 5).Value = 9
$ws.Cells.Item(68, 6).Value = 100112022
$ws.Cells.Item(68, 7).Value = "Arveja Verde"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 45
$ws.Cells.Item(68, 11).Value = 30000
$ws.Cells.Item(68, 12).Value = 32000
$ws.Cells.Item(68, 13).Value = 30667
$ws.Cells.Item(68, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(68, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(68, 16).Value = 1227
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = "Hortaliza"
